# Applies the cryptos list price/volume update described in the commit diff.
# (GitHub Actions scheduled refresh of cryptos.xlsx on 2023-09-17.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.701.92"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.633.37"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'217.98"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "'18.95"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.861.47"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "1.625.86"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").Value = "'63.97"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "26.686.79"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'211.04"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").Value = "'4.29"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "'2.34"
$ws.Range("E22").Value = "  -7.23%  "
$ws.Range("D23").Value = "'6.16"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("D25").Value = "'146.53"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "'7.00"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "'15.49"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "'0.0500"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").Value = "'1.19"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("D34").Value = "1.258.01"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "'2.44"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'0.798"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").Value = "'2.15"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").Value = "1.772.20"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("E44").Value = "  -3.70%  "
$ws.Range("D45").Value = "'90.88"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'59.70"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "'0.407"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.45"
$ws.Range("E51").Value = "  -4.32%  "
